# Fruta / hortaliza, semanal
# A new daily price record was inserted into the "Albahaca" data table at
# row 273 (pushing the existing row 273 and everything below it down by
# one row), and the sheet's used-range dimension grows from R350 to R351.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 273; rows 273:350 shift down to 274:351.
$ws.Range("A273:R273").EntireRow.Insert()

# The newly inserted row starts out blank - duplicate the row immediately
# below it (which now holds what used to be row 273) so every column that
# does NOT change keeps the correct value, then overwrite the columns that
# carry the new record's data.
$ws.Range("A274:R274").Copy()
$ws.Range("A273:R273").PasteSpecial()
$excel.CutCopyMode = $false

# New record's values (date, price range, unit, origin, $/Kg, Kg-or-units).
$ws.Range("D273").Value = 44642
$ws.Range("K273").Value = 2500
$ws.Range("L273").Value = 3000
$ws.Range("M273").Value = 2750
$ws.Range("N273").Value = "$/docena de matas"
$ws.Range("O273").Value = "Región Metropolitana"
$ws.Range("P273").Value = 458
$ws.Range("Q273").Value = 6
